# Correct the 'area_buff' (column E) values for rows 2-101
# on the active worksheet, per the source-data fix described in the
# commit "correct area within buffer".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    277820502.38529497,
    273494554.09821701,
    230965278.75506699,
    249201930.87371299,
    272865554.88393003,
    268613314.09772801,
    249618173.0447,
    267337499.38055599,
    266119167.55680701,
    271903710.572963,
    244809224.61243099,
    274993205.52189898,
    250943925.509101,
    275109655.52770698,
    270722510.95501202,
    264884210.35786599,
    277772517.02549201,
    256238204.68355501,
    272643535.37198001,
    275980406.80915898,
    283191531.606031,
    270948770.57283902,
    225019868.41346499,
    284005273.65626299,
    267437123.77355,
    266387678.25963399,
    254134784.955841,
    281138266.067334,
    268426749.42981699,
    265447383.65678701,
    266028460.27656299,
    261509383.98961699,
    271777503.94406402,
    272547928.83913702,
    279997203.19781899,
    276375031.585356,
    257872995.62318599,
    276201215.574489,
    240332377.002736,
    241512615.73377499,
    280450199.341178,
    278608198.977817,
    268640053.79614198,
    256239805.56469801,
    271251155.47545302,
    254063566.26418799,
    250195467.37032801,
    272026345.142169,
    266543962.58041,
    268865909.50539702,
    261148185.99293599,
    271423781.55661201,
    276682501.89549601,
    269446570.54705298,
    278358563.54608101,
    272667927.83013999,
    252194676.64370999,
    273484936.44867301,
    272491224.92383498,
    280818715.51264399,
    283814453.46674502,
    249634684.14338499,
    263721665.09923601,
    271030642.57695502,
    258780702.277354,
    274910674.203426,
    264548065.92526501,
    248331729.01462299,
    271183268.32600302,
    274795567.11929202,
    246434199.36170101,
    283696727.52098101,
    218850157.463824,
    263765823.37341699,
    278381642.49735802,
    281120746.78577697,
    268948758.23511398,
    274059780.601524,
    281698385.63408399,
    279704335.65184402,
    265187394.163831,
    268592925.05922699,
    269456719.120242,
    277713562.46017802,
    247866618.027392,
    274918501.42144603,
    264622819.21304601,
    277930678.85812098,
    248149684.21466899,
    271505781.30695301,
    280622880.61163598,
    278357946.89252502,
    260480196.44552499,
    280543040.18167001,
    284346056.07838601,
    244451053.14339501,
    277453028.01438498,
    270555111.76331502,
    267706243.112409,
    249198687.217141
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $newValues[$i]
}

# Leave the workbook with the same active-cell selection captured in the
# saved file (M19 on the data sheet).
$ws.Range("M19").Select()
